$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 66.26000000000001
$ws.Range("I15").Value = 66.26000000000001
$ws.Range("K15").Value = 198.78
$ws.Range("M15").Value = -29.78000000000003
$ws.Range("H17").Value = 230553.17
$ws.Range("J17").Value = 230553.17
$ws.Range("L17").Value = 691659.51
$ws.Range("N17").Value = -691995.51
$ws.Range("H32").Value = 699.8182
$ws.Range("J32").Value = 808
$ws.Range("L32").Value = 808
$ws.Range("N32").Value = -1460
$ws.Range("H40").Value = 1609.25
$ws.Range("I40").Value = 1098.4615
$ws.Range("J40").Value = 1855.1852
$ws.Range("K40").Value = 1098.4615
$ws.Range("L40").Value = 1855.1852
$ws.Range("M40").Value = -923.4614999999999
$ws.Range("N40").Value = -2205.1852
$ws.Range("H111").Value = 1133.625
$ws.Range("I111").Value = 939.5
$ws.Range("J111").Value = 1716
$ws.Range("K111").Value = 2818.5
$ws.Range("L111").Value = 5148
$ws.Range("M111").Value = 248.5
$ws.Range("N111").Value = -11282
$ws.Range("H116").Value = 2615.4443
$ws.Range("I116").Value = 1948.4286
$ws.Range("J116").Value = 4950
$ws.Range("K116").Value = 1948.4286
$ws.Range("L116").Value = 4950
$ws.Range("M116").Value = 1493.5714
$ws.Range("N116").Value = -11834
$ws.Range("H125").Value = 1121.4286
$ws.Range("I125").Value = 1020
$ws.Range("J125").Value = 1375
$ws.Range("K125").Value = 9180
$ws.Range("L125").Value = 12375
$ws.Range("M125").Value = -6720
$ws.Range("N125").Value = -17295
$ws.Range("H132").Value = 1274.0695
$ws.Range("I132").Value = 1320.8281
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 3962.4843
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -1432.4843
$ws.Range("N132").Value = -7760
$ws.Range("H137").Value = 964.4286
$ws.Range("I137").Value = 977.7778
$ws.Range("J137").Value = 954.4167
$ws.Range("K137").Value = 2933.3334
$ws.Range("L137").Value = 2863.2501
$ws.Range("M137").Value = -383.3334
$ws.Range("N137").Value = -7963.2501
$ws.Range("H138").Value = 2069.57
$ws.Range("J138").Value = 2638.7273
$ws.Range("L138").Value = 7916.1819
$ws.Range("N138").Value = -18196.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33552.098
$ws.Range("I2").Value = 39724.42
$ws.Range("J2").Value = 1456
$ws.Range("K2").Value = 39724.42
$ws.Range("L2").Value = 1456
$ws.Range("M2").Value = -39611.42
$ws.Range("N2").Value = -1682
$ws.Range("H32").Value = 20190.684
$ws.Range("I32").Value = 19137.5
$ws.Range("J32").Value = 24871.5
$ws.Range("K32").Value = 19137.5
$ws.Range("L32").Value = 24871.5
$ws.Range("M32").Value = -18850.5
$ws.Range("N32").Value = -25445.5
$ws.Range("H61").Value = 803.6875
$ws.Range("I61").Value = 750.97675
$ws.Range("J61").Value = 1257
$ws.Range("K61").Value = 750.97675
$ws.Range("L61").Value = 1257
$ws.Range("M61").Value = -538.97675
$ws.Range("N61").Value = -1681
$ws.Range("H116").Value = 33552.098
$ws.Range("I116").Value = 39724.42
$ws.Range("J116").Value = 1456
$ws.Range("K116").Value = 39724.42
$ws.Range("L116").Value = 1456
$ws.Range("M116").Value = -37430.42
$ws.Range("N116").Value = -6044
$ws.Range("H132").Value = 1817.0571
$ws.Range("I132").Value = 1662.4375
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 4987.3125
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -2457.3125
$ws.Range("N132").Value = -15458.9999
$ws.Range("H136").Value = 803.6875
$ws.Range("I136").Value = 750.97675
$ws.Range("J136").Value = 1257
$ws.Range("K136").Value = 2252.93025
$ws.Range("L136").Value = 3771
$ws.Range("M136").Value = 297.0697499999997
$ws.Range("N136").Value = -8871

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33552.098
$ws.Range("I3").Value = 39724.42
$ws.Range("J3").Value = 1456
$ws.Range("K3").Value = 39724.42
$ws.Range("L3").Value = 1456
$ws.Range("M3").Value = -39610.42
$ws.Range("N3").Value = -1684
$ws.Range("H105").Value = 3144.724
$ws.Range("I105").Value = 3527
$ws.Range("J105").Value = 2295.2222
$ws.Range("K105").Value = 3527
$ws.Range("L105").Value = 2295.2222
$ws.Range("M105").Value = -1780
$ws.Range("N105").Value = -5789.2222
$ws.Range("H134").Value = 16882.303
$ws.Range("I134").Value = 1484.9272
$ws.Range("J134").Value = 93869.17999999999
$ws.Range("K134").Value = 4454.7816
$ws.Range("L134").Value = 281607.54
$ws.Range("M134").Value = -1919.7816
$ws.Range("N134").Value = -286677.54

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1917289.5
$ws.Range("I31").Value = 3002628
$ws.Range("J31").Value = 1986.1765
$ws.Range("K31").Value = 3002628
$ws.Range("L31").Value = 1986.1765
$ws.Range("M31").Value = -3002333
$ws.Range("N31").Value = -2576.1765
$ws.Range("H34").Value = 1917289.5
$ws.Range("I34").Value = 3002628
$ws.Range("J34").Value = 1986.1765
$ws.Range("K34").Value = 3002628
$ws.Range("L34").Value = 1986.1765
$ws.Range("M34").Value = -3002426
$ws.Range("N34").Value = -2390.1765
$ws.Range("H134").Value = 1260.5344
$ws.Range("I134").Value = 1173.738
$ws.Range("K134").Value = 3521.214
$ws.Range("M134").Value = -986.2139999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 137.9375
$ws.Range("I14").Value = 137.9375
$ws.Range("K14").Value = 413.8125
$ws.Range("M14").Value = -240.8125
$ws.Range("H113").Value = 539.9400000000001
$ws.Range("I113").Value = 967.1429000000001
$ws.Range("J113").Value = 470.39536
$ws.Range("K113").Value = 2901.4287
$ws.Range("L113").Value = 1411.18608
$ws.Range("M113").Value = -731.4287000000004
$ws.Range("N113").Value = -5751.186079999999
$ws.Range("H122").Value = 616
$ws.Range("J122").Value = 631.3333
$ws.Range("L122").Value = 5681.9997
$ws.Range("N122").Value = -10581.9997
$ws.Range("H131").Value = 31338.803
$ws.Range("J131").Value = 17337.885
$ws.Range("L131").Value = 52013.655
$ws.Range("N131").Value = -62093.655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3158.6843
$ws.Range("I80").Value = 3100.9375
$ws.Range("J80").Value = 3466.6667
$ws.Range("K80").Value = 3100.9375
$ws.Range("L80").Value = 3466.6667
$ws.Range("M80").Value = -2102.9375
$ws.Range("N80").Value = -5462.6667
$ws.Range("H83").Value = 3158.6843
$ws.Range("I83").Value = 3100.9375
$ws.Range("J83").Value = 3466.6667
$ws.Range("K83").Value = 15504.6875
$ws.Range("L83").Value = 17333.3335
$ws.Range("M83").Value = -10512.6875
$ws.Range("N83").Value = -27317.3335
$ws.Range("H134").Value = 15977.7
$ws.Range("J134").Value = 15977.7
$ws.Range("L134").Value = 47933.10000000001
$ws.Range("N134").Value = -53003.10000000001
$ws.Range("H136").Value = 9383.529
$ws.Range("J136").Value = 9383.529
$ws.Range("L136").Value = 28150.587
$ws.Range("N136").Value = -33250.587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5959.36
$ws.Range("I16").Value = 9070.857
$ws.Range("K16").Value = 9070.857
$ws.Range("M16").Value = -8900.857
$ws.Range("H46").Value = 1888.9286
$ws.Range("I46").Value = 908.8889
$ws.Range("J46").Value = 2353.158
$ws.Range("K46").Value = 908.8889
$ws.Range("L46").Value = 2353.158
$ws.Range("M46").Value = -720.8889
$ws.Range("N46").Value = -2729.158
$ws.Range("H68").Value = 2057
$ws.Range("I68").Value = 2183.1428
$ws.Range("J68").Value = 1762.6666
$ws.Range("K68").Value = 2183.1428
$ws.Range("L68").Value = 1762.6666
$ws.Range("M68").Value = -1434.1428
$ws.Range("N68").Value = -3260.6666
$ws.Range("H71").Value = 2057
$ws.Range("I71").Value = 2183.1428
$ws.Range("J71").Value = 1762.6666
$ws.Range("K71").Value = 10915.714
$ws.Range("L71").Value = 8813.333000000001
$ws.Range("M71").Value = -7171.714
$ws.Range("N71").Value = -16301.333
$ws.Range("H82").Value = 1108.3846
$ws.Range("I82").Value = 969.9231
$ws.Range("J82").Value = 1246.8462
$ws.Range("K82").Value = 969.9231
$ws.Range("L82").Value = 1246.8462
$ws.Range("M82").Value = -608.9231
$ws.Range("N82").Value = -1968.8462
$ws.Range("H85").Value = 1108.3846
$ws.Range("I85").Value = 969.9231
$ws.Range("J85").Value = 1246.8462
$ws.Range("K85").Value = 969.9231
$ws.Range("L85").Value = 1246.8462
$ws.Range("M85").Value = 278.0769
$ws.Range("N85").Value = -3742.8462
$ws.Range("H132").Value = 1860.295
$ws.Range("I132").Value = 1670.8298
$ws.Range("J132").Value = 2496.3572
$ws.Range("K132").Value = 5012.4894
$ws.Range("L132").Value = 7489.071599999999
$ws.Range("M132").Value = -2482.4894
$ws.Range("N132").Value = -12549.0716
$ws.Range("H136").Value = 2092.1592
$ws.Range("I136").Value = 1111.9736
$ws.Range("J136").Value = 8300
$ws.Range("K136").Value = 3335.9208
$ws.Range("L136").Value = 24900
$ws.Range("M136").Value = -785.9207999999999
$ws.Range("N136").Value = -30000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5040.8
$ws.Range("I62").Value = 3001
$ws.Range("J62").Value = 5550.75
$ws.Range("K62").Value = 3001
$ws.Range("L62").Value = 5550.75
$ws.Range("M62").Value = -2377
$ws.Range("N62").Value = -6798.75
$ws.Range("H65").Value = 5040.8
$ws.Range("I65").Value = 3001
$ws.Range("J65").Value = 5550.75
$ws.Range("K65").Value = 15005
$ws.Range("L65").Value = 27753.75
$ws.Range("M65").Value = -11885
$ws.Range("N65").Value = -33993.75
$ws.Range("H110").Value = 28000
$ws.Range("J110").Value = 28000
$ws.Range("L110").Value = 28000
$ws.Range("N110").Value = -36180
$ws.Range("H132").Value = 957.9697
$ws.Range("I132").Value = 750.1539
$ws.Range("J132").Value = 1729.8572
$ws.Range("K132").Value = 2250.4617
$ws.Range("L132").Value = 5189.571599999999
$ws.Range("M132").Value = 279.5383000000002
$ws.Range("N132").Value = -10249.5716
$ws.Range("H136").Value = 1831.9333
$ws.Range("I136").Value = 1891.3572
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 5674.071599999999
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -3124.071599999999
$ws.Range("N136").Value = -8100
$ws.Range("H140").Value = 53911.152
$ws.Range("J140").Value = 53911.152
$ws.Range("L140").Value = 53911.152
$ws.Range("N140").Value = -64271.152
